## Adds a new paragraph ("B.\tDit is de <del>tekst</del>") right after the
## existing empty paragraph, with "tekst" recorded as a tracked deletion
## (w:del / w:delText) authored by thijs@label305.onmicrosoft.com.

$d = $word.ActiveDocument

# The document has exactly two paragraphs to start: the "A. ..." paragraph
# and a trailing empty paragraph. Locate the trailing empty paragraph
# (last paragraph) robustly instead of hard-coding an index.
$lastPara = $d.Paragraphs.Last

# Append a brand-new paragraph right after it; Word seeds it with the same
# paragraph formatting (pPr/rPr) as the paragraph it was split from.
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last

# Build the whole new paragraph's text in one go with plain COM calls so it
# ends up as a single contiguous range (avoids leaving stray empty runs
# behind). The tracked deletion text ("tekst") is appended straight onto
# the end of "Dit is de ".
$body = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$body.InsertAfter("B." + [char]9 + "Dit is de tekst")

# Re-select that same (now non-empty) paragraph content range and replace
# it wholesale with the precise run/tab/tracked-deletion structure we want:
#   <w:r>B.</w:r><w:r><w:tab/>Dit is de </w:r><w:del>...tekst...</w:del>
$full = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:rPr>
                <w:szCs w:val="22"/>
              </w:rPr>
              <w:t>B.</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:szCs w:val="22"/>
              </w:rPr>
              <w:tab/>
              <w:t xml:space="preserve">Dit is de </w:t>
            </w:r>
            <w:del w:id="3" w:author="thijs@label305.onmicrosoft.com" w:date="2021-02-26T14:25:00Z">
              <w:r>
                <w:rPr>
                  <w:szCs w:val="22"/>
                </w:rPr>
                <w:delText>tekst</w:delText>
              </w:r>
            </w:del>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$full.InsertXML($xml)
